# Renumber the "Id" column (A) on the first sheet from the legacy
# 724-733 sequence down to a simple 1-9 counter (row 11 becomes 90),
# then move the viewport/selection up to reflect where editing resumed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$newValues = @{
    2  = 1
    3  = 2
    4  = 3
    5  = 4
    6  = 5
    7  = 6
    8  = 7
    9  = 8
    10 = 9
    11 = 90
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 1).Value = $newValues[$row]
}

# Update the sheet view: scrolled-to row and current selection.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("B13").Select()
